$wb = $excel.ActiveWorkbook

# Update the "Status" value from "Ready for handoff" to "In Translation"
# on all three sheets (Overview, zh-cn, de-de), then autofit the affected
# columns so their widths recompute to match the shorter text.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Range("E:F").EntireColumn.AutoFit() | Out-Null

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Range("C:C").EntireColumn.AutoFit() | Out-Null

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Range("C:C").EntireColumn.AutoFit() | Out-Null
